$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B (old "Betrag" column shifts from B to C)
$ws.Columns("B:B").Insert()

# The insert inherited column A's date style (s=2) - clear that back to default
$ws.Range("B1:B9").ClearFormats()

# New column header + Belegnr. (invoice number) values
$ws.Range("B1").Value = "Belegnr."
$ws.Range("B2").Value = 6000001
$ws.Range("B3").Value = 6000002
$ws.Range("B4").Value = 6000003
$ws.Range("B5").Value = 6000004
$ws.Range("B6").Value = 6000008
$ws.Range("B7").Value = 6000006
$ws.Range("B8").Value = 6000010
$ws.Range("B9").Value = 6000011

# Fix the Betrag value that changed precision in row 6 (now column C)
$ws.Range("C6").Value = -390.05

# New row 9 of data - clone formatting from row 8 so the same cell styles are reused
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = [DateTime]::FromOADate(43160)

$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = -250

# Column C (Betrag) formatting/width, matching the numeric style used elsewhere in that column
$ws.Columns("C:C").ColumnWidth = 11

# Stray formatted cells left over from selecting a wide range while formatting
$ws.Range("I1:I2").NumberFormat = $ws.Range("C2").NumberFormat

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the active selection shown in the final workbook
$ws.Range("C7").Select()
